# Generate Report for Handback
# Update status cells from "Ready for handoff" to "Handed back: in sync with en-US"
# and refresh the "Latest Handback DateTime" values for the handed-back files.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: row 3 corresponds to e27be0ea-...md, mark it handed back.
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# zh-cn detail sheet: update status and handback datetime for rows 2 and 3.
$zhcn.Range("B3").Value = $handedBack
$zhcn.Range("G2").Value = "2016-03-04 04:11:42"
$zhcn.Range("G3").Value = "2016-03-04 04:11:42"

# de-de detail sheet: update status and handback datetime for rows 2 and 3.
$dede.Range("B3").Value = $handedBack
$dede.Range("G2").Value = "2016-03-04 04:12:12"
$dede.Range("G3").Value = "2016-03-04 04:12:12"
